$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 24,14
$values[0,0] = 0.3770835877108993
$values[0,1] = 0.06204085455566144
$values[0,2] = 0.02428201727138912
$values[0,3] = 0.4169131172345146
$values[0,4] = 0.674504049928629
$values[0,5] = 0
$values[0,6] = 0.07973214163530429
$values[0,7] = 0
$values[0,8] = 0
$values[0,9] = 0.362824689027974
$values[0,10] = 0
$values[0,11] = 0
$values[0,12] = 1.274553151604099
$values[0,13] = 2.286437771945032
$values[1,0] = 0.3350854292761483
$values[1,1] = 0.05623166161785775
$values[1,2] = 0.02260588466011626
$values[1,3] = 0.3637838837403109
$values[1,4] = 0.6703514289716281
$values[1,5] = 0
$values[1,6] = 0.07973214163530429
$values[1,7] = 0
$values[1,8] = 0
$values[1,9] = 0.3177370021360844
$values[1,10] = 0
$values[1,11] = 0
$values[1,12] = 1.290348779400875
$values[1,13] = 2.286328194739156
$values[2,0] = 0.3093405411897834
$values[2,1] = 0.05264010522283513
$values[2,2] = 0.02156747538315784
$values[2,3] = 0.3312461650622396
$values[2,4] = 0.6682533619767952
$values[2,5] = 0
$values[2,6] = 0.07973214163530429
$values[2,7] = 0
$values[2,8] = 0
$values[2,9] = 0.2900296133697111
$values[2,10] = 0
$values[2,11] = 0
$values[2,12] = 1.300534452545329
$values[2,13] = 2.287740803486514
$values[3,0] = 0.2988602305306642
$values[3,1] = 0.05117033785384706
$values[3,2] = 0.02114201639177082
$values[3,3] = 0.3180062801571495
$values[3,4] = 0.6675118070990322
$values[3,5] = 0
$values[3,6] = 0.07973214163530429
$values[3,7] = 0
$values[3,8] = 0
$values[3,9] = 0.2787330757165876
$values[3,10] = 0
$values[3,11] = 0
$values[3,12] = 1.304807596311655
$values[3,13] = 2.28868799294176
$values[4,0] = 0.297120653797009
$values[4,1] = 0.05092591177715633
$values[4,2] = 0.02107123126454979
$values[4,3] = 0.3158089389432064
$values[4,4] = 0.6673955187060301
$values[4,5] = 0
$values[4,6] = 0.07973214163530429
$values[4,7] = 0
$values[4,8] = 0
$values[4,9] = 0.2768569685403293
$values[4,10] = 0
$values[4,11] = 0
$values[4,12] = 1.305524538768959
$values[4,13] = 2.28886769594223
$values[5,0] = 0.3091991553052935
$values[5,1] = 0.0526203084050394
$values[5,2] = 0.02156174676378697
$values[5,3] = 0.3310675304697952
$values[5,4] = 0.6682429020566616
$values[5,5] = 0
$values[5,6] = 0.07973214163530429
$values[5,7] = 0
$values[5,8] = 0
$values[5,9] = 0.2898772863597685
$values[5,10] = 0
$values[5,11] = 0
$values[5,12] = 1.300591586181338
$values[5,13] = 2.287752074055049
$values[6,0] = 0.3625940654347346
$values[6,1] = 0.06004298466099556
$values[6,2] = 0.02370602493705576
$values[6,3] = 0.3985757854852636
$values[6,4] = 0.6729783912405267
$values[6,5] = 0
$values[6,6] = 0.07973214163530429
$values[6,7] = 0
$values[6,8] = 0
$values[6,9] = 0.347283484634147
$values[6,10] = 0
$values[6,11] = 0
$values[6,12] = 1.279898338295507
$values[6,13] = 2.286092475541409
$values[7,0] = 0.4676261129159798
$values[7,1] = 0.07440233283342934
$values[7,2] = 0.02783640828050693
$values[7,3] = 0.5317021362696437
$values[7,4] = 0.6858566540029898
$values[7,5] = 0
$values[7,6] = 0.07973214163530429
$values[7,7] = 0
$values[7,8] = 0
$values[7,9] = 0.459662796498634
$values[7,10] = 0
$values[7,11] = 0
$values[7,12] = 1.243186959970292
$values[7,13] = 2.294610304080578
$values[8,0] = 0.5449860507108895
$values[8,1] = 0.08483253962364756
$values[8,2] = 0.03082436763538965
$values[8,3] = 0.6300792176655534
$values[8,4] = 0.697521776632442
$values[8,5] = 0
$values[8,6] = 0.07973214163530429
$values[8,7] = 0
$values[8,8] = 0
$values[8,9] = 0.5421061580285027
$values[8,10] = 0
$values[8,11] = 0
$values[8,12] = 1.218576562687081
$values[8,13] = 2.308091560551418
$values[9,0] = 0.5802203928164431
$values[9,1] = 0.0895516444762734
$values[9,2] = 0.03217330660627482
$values[9,3] = 0.6749827859718778
$values[9,4] = 0.7033101054859827
$values[9,5] = 0
$values[9,6] = 0.07973214163530429
$values[9,7] = 0
$values[9,8] = 0
$values[9,9] = 0.5795851526362696
$values[9,10] = 0
$values[9,11] = 0
$values[9,12] = 1.207894173695783
$values[9,13] = 2.3158033826561
$values[10,0] = 0.5935686640659981
$values[10,1] = 0.09133494437777756
$values[10,2] = 0.032682607991255
$values[10,3] = 0.692010265313499
$values[10,4] = 0.7055714866965985
$values[10,5] = 0
$values[10,6] = 0.07973214163530429
$values[10,7] = 0
$values[10,8] = 0
$values[10,9] = 0.5937736860145719
$values[10,10] = 0
$values[10,11] = 0
$values[10,12] = 1.203922902070837
$values[10,13] = 2.318951489422119
$values[11,0] = 0.5906936235362821
$values[11,1] = 0.09095104494875272
$values[11,2] = 0.03257298859691105
$values[11,3] = 0.688342021650584
$values[11,4] = 0.7050813653103916
$values[11,5] = 0
$values[11,6] = 0.07973214163530429
$values[11,7] = 0
$values[11,8] = 0
$values[11,9] = 0.5907181129259129
$values[11,10] = 0
$values[11,11] = 0
$values[11,12] = 1.204774896588745
$values[11,13] = 2.318263344849242
$values[12,0] = 0.5813184492519952
$values[12,1] = 0.08969843243542641
$values[12,2] = 0.03221523763128431
$values[12,3] = 0.6763831666736309
$values[12,4] = 0.7034947574093309
$values[12,5] = 0
$values[12,6] = 0.07973214163530429
$values[12,7] = 0
$values[12,8] = 0
$values[12,9] = 0.5807525330477006
$values[12,10] = 0
$values[12,11] = 0
$values[12,12] = 1.207565972518413
$values[12,13] = 2.316057809504969
$values[13,0] = 0.5755766271189202
$values[13,1] = 0.0889306852835432
$values[13,2] = 0.03199590692432963
$values[13,3] = 0.6690611305115368
$values[13,4] = 0.7025319663316054
$values[13,5] = 0
$values[13,6] = 0.07973214163530429
$values[13,7] = 0
$values[13,8] = 0
$values[13,9] = 0.574647802540909
$values[13,10] = 0
$values[13,11] = 0
$values[13,12] = 1.209285217770342
$values[13,13] = 2.314736544367577
$values[14,0] = 0.5426842346172407
$values[14,1] = 0.08452361595524849
$values[14,2] = 0.03073600136175259
$values[14,3] = 0.6271478550841465
$values[14,4] = 0.6971532065754076
$values[14,5] = 0
$values[14,6] = 0.07973214163530429
$values[14,7] = 0
$values[14,8] = 0
$values[14,9] = 0.539656285018566
$values[14,10] = 0
$values[14,11] = 0
$values[14,12] = 1.219285017447427
$values[14,13] = 2.307619413299136
$values[15,0] = 0.522516538096653
$values[15,1] = 0.08181342782572187
$values[15,2] = 0.0299604298429017
$values[15,3] = 0.6014753923363543
$values[15,4] = 0.693977044181679
$values[15,5] = 0
$values[15,6] = 0.07973214163530429
$values[15,7] = 0
$values[15,8] = 0
$values[15,9] = 0.5181834819954645
$values[15,10] = 0
$values[15,11] = 0
$values[15,12] = 1.225551067461026
$values[15,13] = 2.30365825535074
$values[16,0] = 0.510920674849018
$values[16,1] = 0.0802521871741817
$values[16,2] = 0.02951337413139044
$values[16,3] = 0.5867234126350525
$values[16,4] = 0.69219553842143
$values[16,5] = 0
$values[16,6] = 0.07973214163530429
$values[16,7] = 0
$values[16,8] = 0
$values[16,9] = 0.5058305713976665
$values[16,10] = 0
$values[16,11] = 0
$values[16,12] = 1.229203411001974
$values[16,13] = 2.301528493406948
$values[17,0] = 0.5069952241472606
$values[17,1] = 0.07972316468554652
$values[17,2] = 0.02936184354677351
$values[17,3] = 0.5817310227705548
$values[17,4] = 0.6916001325926686
$values[17,5] = 0
$values[17,6] = 0.07973214163530429
$values[17,7] = 0
$values[17,8] = 0
$values[17,9] = 0.5016477004940327
$values[17,10] = 0
$values[17,11] = 0
$values[17,12] = 1.230448318999632
$values[17,13] = 2.30083289115251
$values[18,0] = 0.5246630061439816
$values[18,1] = 0.08210218187666385
$values[18,2] = 0.03004309111527448
$values[18,3] = 0.6042067938678031
$values[18,4] = 0.6943104583603912
$values[18,5] = 0
$values[18,6] = 0.07973214163530429
$values[18,7] = 0
$values[18,8] = 0
$values[18,9] = 0.5204695415436333
$values[18,10] = 0
$values[18,11] = 0
$values[18,12] = 1.224879038154809
$values[18,13] = 2.304064543444042
$values[19,0] = 0.5840720110724078
$values[19,1] = 0.09006645640147326
$values[19,2] = 0.03232035906292197
$values[19,3] = 0.6798951215308051
$values[19,4] = 0.7039588958371326
$values[19,5] = 0
$values[19,6] = 0.07973214163530429
$values[19,7] = 0
$values[19,8] = 0
$values[19,9] = 0.5836797755066243
$values[19,10] = 0
$values[19,11] = 0
$values[19,12] = 1.206744158572588
$values[19,13] = 2.316699440240171
$values[20,0] = 0.6229328061783121
$values[20,1] = 0.09524985580731027
$values[20,2] = 0.0337998545811331
$values[20,3] = 0.7294994389327201
$values[20,4] = 0.7106696982456526
$values[20,5] = 0
$values[20,6] = 0.07973214163530429
$values[20,7] = 0
$values[20,8] = 0
$values[20,9] = 0.6249681334513184
$values[20,10] = 0
$values[20,11] = 0
$values[20,12] = 1.195322883767823
$values[20,13] = 2.326285251598364
$values[21,0] = 0.6021891218030362
$values[21,1] = 0.09248537614878671
$values[21,2] = 0.03301103839224595
$values[21,3] = 0.7030115201910832
$values[21,4] = 0.7070508991911169
$values[21,5] = 0
$values[21,6] = 0.07973214163530429
$values[21,7] = 0
$values[21,8] = 0
$values[21,9] = 0.6029340020307359
$values[21,10] = 0
$values[21,11] = 0
$values[21,12] = 1.201379154373702
$values[21,13] = 2.321047353218432
$values[22,0] = 0.5236925916054531
$values[22,1] = 0.08197164587932093
$values[22,2] = 0.03000572359363929
$values[22,3] = 0.6029719041564476
$values[22,4] = 0.6941595831670782
$values[22,5] = 0
$values[22,6] = 0.07973214163530429
$values[22,7] = 0
$values[22,8] = 0
$values[22,9] = 0.5194360386354333
$values[22,10] = 0
$values[22,11] = 0
$values[22,12] = 1.225182707146159
$values[22,13] = 2.303880400991517
$values[23,0] = 0.4391779234190381
$values[23,1] = 0.07053876740212672
$values[23,2] = 0.0267271338740116
$values[23,3] = 0.4955964288448769
$values[23,4] = 0.6819867835137501
$values[23,5] = 0
$values[23,6] = 0.07973214163530429
$values[23,7] = 0
$values[23,8] = 0
$values[23,9] = 0.4292823501930059
$values[23,10] = 0
$values[23,11] = 0
$values[23,12] = 1.252703924108705
$values[23,13] = 2.291040909878177

$ws.Range("B2:O25").Value = $values

Write-Host "Done writing values"